# Daily attendance processing - 2026-01-15 22:36:52
# Reorders the "Recorded By" (column G) names for specific session rows so
# that "dnasr281@gmail.com" is listed second instead of first, e.g.
#   "dnasr281@gmail.com, System"        -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $current = [string]$cell.Value2

    $parts = $current -split ", ", 2

    if ($parts.Count -eq 2 -and $parts[0] -eq "dnasr281@gmail.com") {
        $newValue = $parts[1] + ", " + $parts[0]
        $cell.Value = $newValue
    }
}
